$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet - refreshed loan numbers (rounding changed slightly after the
# automation script was re-run / stabilized)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 1402.28
$wsSummary.Range("E2").Value = 8597.7199999999993
$wsSummary.Range("F2").Value = 276.18

$wsSummary.Range("A3").Value = 563.33000000000004
$wsSummary.Range("E3").Value = 466.35

$wsSummary.Range("A5").Value = 0.74
$wsSummary.Range("B5").Value = 0.74

# ---------------------------------------------------------------------------
# Repayment schedule sheet - same re-run, updated installment figures
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

$wsSchedule.Range("J5").Value = 0.74
$wsSchedule.Range("K5").Value = 888.46
$wsSchedule.Range("Q5").Value = 276.18

$wsSchedule.Range("F6").Value = 800.09
$wsSchedule.Range("G6").Value = 7521.45
$wsSchedule.Range("H6").Value = 87.63

$wsSchedule.Range("G7").Value = 6707.91

$wsSchedule.Range("G8").Value = 5888.56

$wsSchedule.Range("G9").Value = 5058.92

$wsSchedule.Range("G10").Value = 4222.76

$wsSchedule.Range("G11").Value = 3378.08

$wsSchedule.Range("G12").Value = 2523.6799999999998

$wsSchedule.Range("G13").Value = 1661.68

$wsSchedule.Range("G14").Value = 790.35

$wsSchedule.Range("F15").Value = 790.35
$wsSchedule.Range("H15").Value = 8.06
$wsSchedule.Range("K15").Value = 798.41
$wsSchedule.Range("Q15").Value = 798.41

# ---------------------------------------------------------------------------
# Transactions sheet - renumbered transaction ids + refreshed amounts
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Range("A2").Value = 446
$wsTransactions.Range("E2").Value = 87.63
$wsTransactions.Range("G2").Value = 87.63

$wsTransactions.Range("A3").Value = 445
$wsTransactions.Range("E3").Value = 23.75
$wsTransactions.Range("I3").Value = 0.74

$wsTransactions.Range("A4").Value = 443
$wsTransactions.Range("J4").Value = 8597.7199999999993

$wsTransactions.Range("A5").Value = 442
$wsTransactions.Range("F5").Value = 1402.28
$wsTransactions.Range("I5").Value = 0.74
$wsTransactions.Range("J5").Value = 3597.72

$wsTransactions.Range("A6").Value = 444

$wsTransactions.Range("A7").Value = 441

# ---------------------------------------------------------------------------
# Restore the selection on each sheet (recorded cursor position at save time)
# and finish with "Transactions" as the active sheet/tab, same as before.
# ---------------------------------------------------------------------------
$wsSummary.Range("C8").Select() | Out-Null
$wsSchedule.Range("L6").Select() | Out-Null
$wsTransactions.Range("D10").Select() | Out-Null
